$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as text (matches source data)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "25.213.03"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "1.550.76"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "206.51"
$ws.Range("E5").Value = "  -3.52%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -5.56%  "
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("D10").Value = "17.77"
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "1.766.19"
$ws.Range("E12").Value = "  -4.81%  "
$ws.Range("D13").Value = "1.547.42"
$ws.Range("E13").Value = "  -11.94%  "
$ws.Range("D14").Value = "3.97"
$ws.Range("E14").Value = "  -5.31%  "
$ws.Range("D15").Value = "0.502"
$ws.Range("E15").Value = "  -4.98%  "
$ws.Range("D16").Value = "25.181.11"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").Value = "0.0₃0707"
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "58.58"
$ws.Range("E18").Value = "  -4.53%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "185.21"
$ws.Range("E20").Value = "  -3.66%  "
$ws.Range("D21").Value = "4.09"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("D26").Value = "139.35"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").Value = "1.64"
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("D28").Value = "14.74"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "6.37"
$ws.Range("E29").Value = "  -5.43%  "
$ws.Range("E30").Value = "  -6.94%  "
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("D34").Value = "1.44"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "1.085.04"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("E39").Value = "  -5.66%  "
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.757"
$ws.Range("E42").Value = "  -11.40%  "
$ws.Range("D43").Value = "92.52"
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("D44").Value = "5.03"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").Value = "1.682.68"
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("E46").Value = "  +11.94%  "
$ws.Range("D47").Value = "52.17"
$ws.Range("E47").Value = "  -4.19%  "
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("E49").Value = "  -5.77%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -1.99%  "
